$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search DropDown")
$ws.Activate()

# Fix the PersonFilter value for Population Health Sciences (was previously
# stored without spaces). The dependent formula in column B recalculates
# automatically since it references A4.
$ws.Range("A4").Value = "Population Health Sciences"

# Leave the sheet scrolled/selected where the edit was made.
$ws.Range("B18").Select()
